$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rows 25-30 currently carry the old "empty placeholder" formatting (s=18/19).
# Copy the formatting used by the already-correct rows 19-24 down onto 25-30
# before writing values, so the new rows look like a continuation of the table.
$ws.Range("A19:Q19").Copy()
$ws.Range("A25:Q30").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Row 19
$ws.Range("A19").Value = "8088"
$ws.Range("B19").Value = "1001"
$ws.Range("C19").Value = "8088"
$ws.Range("D19").Value = "0000"
$ws.Range("E19").Value = "aarch64"
$ws.Range("F19").Value = "openEuler 22.03 LTS SP3"
$ws.Range("G19").Value = "txgbe"
$ws.Range("H19").Value = "1.3.2oe"
$ws.Range("I19").Value = "NIC"
$ws.Range("J19").Value = "2023.12.05"
$ws.Range("K19").Value = "91049901f1c0c1f717b646505ffb9c066fda004932e6031ab5dcc7ed336f9358"
$ws.Range("L19").Value = "95K"
$ws.Range("M19").Value = "Netswift"
$ws.Range("N19").Value = "RP1000P2SFP"
$ws.Range("O19").Value = "SP1000A"
$ws.Range("Q19").Value = "inbox"

# Row 20
$ws.Range("A20").Value = "8088"
$ws.Range("B20").Value = "2001"
$ws.Range("C20").Value = "8088"
$ws.Range("D20").Value = "2000"
$ws.Range("E20").Value = "aarch64"
$ws.Range("F20").Value = "openEuler 22.03 LTS SP3"
$ws.Range("G20").Value = "txgbe"
$ws.Range("H20").Value = "1.3.2oe"
$ws.Range("I20").Value = "NIC"
$ws.Range("J20").Value = "2023.12.05"
$ws.Range("K20").Value = "91049901f1c0c1f717b646505ffb9c066fda004932e6031ab5dcc7ed336f9358"
$ws.Range("L20").Value = "95K"
$ws.Range("M20").Value = "Netswift"
$ws.Range("N20").Value = "RP2000P2SFP"
$ws.Range("O20").Value = "WX1820AL"
$ws.Range("Q20").Value = "inbox"

# Row 21
$ws.Range("A21").Value = "8088"
$ws.Range("B21").Value = "0105"
$ws.Range("C21").Value = "8088"
$ws.Range("D21").Value = "0202"
$ws.Range("E21").Value = "aarch64"
$ws.Range("F21").Value = "openEuler 22.03 LTS SP3"
$ws.Range("G21").Value = "ngbe"
$ws.Range("H21").Value = "1.2.2oe"
$ws.Range("I21").Value = "NIC"
$ws.Range("J21").Value = "2023.12.05"
$ws.Range("K21").Value = "9b3a6c73035fcdd4596c8236c827369a01778ab3ce807a2451c7c98fc4cbe9f8"
$ws.Range("L21").Value = "87K"
$ws.Range("M21").Value = "Netswift"
$ws.Range("N21").Value = "SF200HT"
$ws.Range("O21").Value = "WX1860AL2"
$ws.Range("Q21").Value = "inbox"

# Row 22
$ws.Range("A22").Value = "8088"
$ws.Range("B22").Value = "0101"
$ws.Range("C22").Value = "8088"
$ws.Range("D22").Value = "0201"
$ws.Range("E22").Value = "aarch64"
$ws.Range("F22").Value = "openEuler 22.03 LTS SP3"
$ws.Range("G22").Value = "ngbe"
$ws.Range("H22").Value = "1.2.2oe"
$ws.Range("I22").Value = "NIC"
$ws.Range("J22").Value = "2023.12.05"
$ws.Range("K22").Value = "9b3a6c73035fcdd4596c8236c827369a01778ab3ce807a2451c7c98fc4cbe9f8"
$ws.Range("L22").Value = "87K"
$ws.Range("M22").Value = "Netswift"
$ws.Range("N22").Value = "SF200T"
$ws.Range("O22").Value = "WX1860A2"
$ws.Range("Q22").Value = "inbox"

# Row 23
$ws.Range("A23").Value = "8088"
$ws.Range("B23").Value = "0107"
$ws.Range("C23").Value = "8088"
$ws.Range("D23").Value = "0402"
$ws.Range("E23").Value = "aarch64"
$ws.Range("F23").Value = "openEuler 22.03 LTS SP3"
$ws.Range("G23").Value = "ngbe"
$ws.Range("H23").Value = "1.2.2oe"
$ws.Range("I23").Value = "NIC"
$ws.Range("J23").Value = "2023.12.05"
$ws.Range("K23").Value = "9b3a6c73035fcdd4596c8236c827369a01778ab3ce807a2451c7c98fc4cbe9f8"
$ws.Range("L23").Value = "87K"
$ws.Range("M23").Value = "Netswift"
$ws.Range("N23").Value = "SF400HT"
$ws.Range("O23").Value = "WX1860AL4"
$ws.Range("Q23").Value = "inbox"

# Row 24
$ws.Range("A24").Value = "8088"
$ws.Range("B24").Value = "0107"
$ws.Range("C24").Value = "8088"
$ws.Range("D24").Value = "0401"
$ws.Range("E24").Value = "aarch64"
$ws.Range("F24").Value = "openEuler 22.03 LTS SP3"
$ws.Range("G24").Value = "ngbe"
$ws.Range("H24").Value = "1.2.2oe"
$ws.Range("I24").Value = "NIC"
$ws.Range("J24").Value = "2023.12.05"
$ws.Range("K24").Value = "9b3a6c73035fcdd4596c8236c827369a01778ab3ce807a2451c7c98fc4cbe9f8"
$ws.Range("L24").Value = "87K"
$ws.Range("M24").Value = "Netswift"
$ws.Range("N24").Value = "SF400T"
$ws.Range("O24").Value = "WX1860A4"
$ws.Range("Q24").Value = "inbox"

# Row 25
$ws.Range("A25").Value = "8088"
$ws.Range("B25").Value = "1001"
$ws.Range("C25").Value = "8088"
$ws.Range("D25").Value = "0000"
$ws.Range("E25").Value = "x86_64"
$ws.Range("F25").Value = "openEuler 22.03 LTS SP3"
$ws.Range("G25").Value = "txgbe"
$ws.Range("H25").Value = "1.3.2oe"
$ws.Range("I25").Value = "NIC"
$ws.Range("J25").Value = "2023.12.05"
$ws.Range("K25").Value = "39d480b1ca092fb4bc01732d077525a595b8b2ed3b63fb5d9c6bd21dab290d1b"
$ws.Range("L25").Value = "107K"
$ws.Range("M25").Value = "Netswift"
$ws.Range("N25").Value = "RP1000P2SFP"
$ws.Range("O25").Value = "SP1000A"
$ws.Range("Q25").Value = "inbox"

# Row 26
$ws.Range("A26").Value = "8088"
$ws.Range("B26").Value = "2001"
$ws.Range("C26").Value = "8088"
$ws.Range("D26").Value = "2000"
$ws.Range("E26").Value = "x86_64"
$ws.Range("F26").Value = "openEuler 22.03 LTS SP3"
$ws.Range("G26").Value = "txgbe"
$ws.Range("H26").Value = "1.3.2oe"
$ws.Range("I26").Value = "NIC"
$ws.Range("J26").Value = "2023.12.05"
$ws.Range("K26").Value = "39d480b1ca092fb4bc01732d077525a595b8b2ed3b63fb5d9c6bd21dab290d1b"
$ws.Range("L26").Value = "107K"
$ws.Range("M26").Value = "Netswift"
$ws.Range("N26").Value = "RP2000P2SFP"
$ws.Range("O26").Value = "WX1820AL"
$ws.Range("Q26").Value = "inbox"

# Row 27
$ws.Range("A27").Value = "8088"
$ws.Range("B27").Value = "0105"
$ws.Range("C27").Value = "8088"
$ws.Range("D27").Value = "0202"
$ws.Range("E27").Value = "x86_64"
$ws.Range("F27").Value = "openEuler 22.03 LTS SP3"
$ws.Range("G27").Value = "ngbe"
$ws.Range("H27").Value = "1.2.2oe"
$ws.Range("I27").Value = "NIC"
$ws.Range("J27").Value = "2023.12.05"
$ws.Range("K27").Value = "b0c9650549a49fc6bc1b37a6a1f06a528f41d13174b2777883a0da03ab435fc4"
$ws.Range("L27").Value = "97K"
$ws.Range("M27").Value = "Netswift"
$ws.Range("N27").Value = "SF200HT"
$ws.Range("O27").Value = "WX1860AL2"
$ws.Range("Q27").Value = "inbox"

# Row 28
$ws.Range("A28").Value = "8088"
$ws.Range("B28").Value = "0101"
$ws.Range("C28").Value = "8088"
$ws.Range("D28").Value = "0201"
$ws.Range("E28").Value = "x86_64"
$ws.Range("F28").Value = "openEuler 22.03 LTS SP3"
$ws.Range("G28").Value = "ngbe"
$ws.Range("H28").Value = "1.2.2oe"
$ws.Range("I28").Value = "NIC"
$ws.Range("J28").Value = "2023.12.05"
$ws.Range("K28").Value = "b0c9650549a49fc6bc1b37a6a1f06a528f41d13174b2777883a0da03ab435fc4"
$ws.Range("L28").Value = "97K"
$ws.Range("M28").Value = "Netswift"
$ws.Range("N28").Value = "SF200T"
$ws.Range("O28").Value = "WX1860A2"
$ws.Range("Q28").Value = "inbox"

# Row 29
$ws.Range("A29").Value = "8088"
$ws.Range("B29").Value = "0107"
$ws.Range("C29").Value = "8088"
$ws.Range("D29").Value = "0402"
$ws.Range("E29").Value = "x86_64"
$ws.Range("F29").Value = "openEuler 22.03 LTS SP3"
$ws.Range("G29").Value = "ngbe"
$ws.Range("H29").Value = "1.2.2oe"
$ws.Range("I29").Value = "NIC"
$ws.Range("J29").Value = "2023.12.05"
$ws.Range("K29").Value = "b0c9650549a49fc6bc1b37a6a1f06a528f41d13174b2777883a0da03ab435fc4"
$ws.Range("L29").Value = "97K"
$ws.Range("M29").Value = "Netswift"
$ws.Range("N29").Value = "SF400HT"
$ws.Range("O29").Value = "WX1860AL4"
$ws.Range("Q29").Value = "inbox"

# Row 30
$ws.Range("A30").Value = "8088"
$ws.Range("B30").Value = "0103"
$ws.Range("C30").Value = "8088"
$ws.Range("D30").Value = "0401"
$ws.Range("E30").Value = "x86_64"
$ws.Range("F30").Value = "openEuler 22.03 LTS SP3"
$ws.Range("G30").Value = "ngbe"
$ws.Range("H30").Value = "1.2.2oe"
$ws.Range("I30").Value = "NIC"
$ws.Range("J30").Value = "2023.12.05"
$ws.Range("K30").Value = "b0c9650549a49fc6bc1b37a6a1f06a528f41d13174b2777883a0da03ab435fc4"
$ws.Range("L30").Value = "97K"
$ws.Range("M30").Value = "Netswift"
$ws.Range("N30").Value = "SF400T"
$ws.Range("O30").Value = "WX1860A4"
$ws.Range("Q30").Value = "inbox"

$ws.Activate()
$ws.Range("N33").Select()
